$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: boolean "published" style flags for each course row
$ws.Range("D1").Value = $false
$ws.Range("D2").Value = $false
$ws.Range("D3").Value = $true
$ws.Range("D4").Value = $true

# New column F: an empty cell that only carries right-aligned formatting
$ws.Range("F2").HorizontalAlignment = -4152

# Move the active selection to E2
$ws.Range("E2").Select() | Out-Null
